$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 9 new data rows (22-30), continuing the existing table pattern:
# A = regcntr_id, B = machine_id, C = lang_code ("eng"), D = is_active (TRUE),
# E = cr_by ("superadmin"), F = cr_dtimes ("now()"), G = eff_dtimes ("now()")
$regIds  = @(10002, 10003, 10004, 10005, 10006, 10007, 10008, 10009, 10010)
$machIds = @(10021, 10022, 10023, 10024, 10025, 10026, 10027, 10028, 10029)

for ($i = 0; $i -lt 9; $i++) {
    $r = 22 + $i
    $ws.Cells.Item($r, 1).Value = $regIds[$i]
    $ws.Cells.Item($r, 2).Value = $machIds[$i]
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Update the sheet view: scroll the window so row 19 is at the top and
# select the newly-added machine_id column (B22:B30, active cell B22).
$ws.Range("B22:B30").Select()
$excel.ActiveWindow.ScrollRow = 19

# Mark the print orientation as portrait so the worksheet gets a pageSetup entry.
$ws.PageSetup.Orientation = 1
